# Edit tbl8 (sheet8) per the target diff:
#  - drop the old row 7 ("ss6" strategy-count row) entirely, shifting the
#    aic/bic/pr_chisq summary rows up
#  - rename the "0.5_*" columns to "final_*"
#  - refresh the numeric coefficient/p-value strings for the remaining rows
#  - clear the now-unused "final" (G/H) column values and the trailing
#    pr_chisq row so the used range shrinks to A1:H8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tbl8")

# Remove the old row 7 (strategy "ss6"); rows 8-10 shift up to 7-9.
$ws.Rows(7).Delete()

# Header renames
$ws.Range("G1").Value2 = "final_Coefficient (95% CI)"
$ws.Range("H1").Value2 = "final_p-value"

# Row 2 (ss1)
$ws.Range("C2").Value2 = "0.15 (-0.38, 0.69)"
$ws.Range("E2").Value2 = "1.06 (0.34, 1.84)"
$ws.Range("G2:H2").ClearContents()

# Row 3 (ss2)
$ws.Range("C3").Value2 = "-0.23 (-0.79, 0.33)"
$ws.Range("E3").Value2 = "-0.2 (-0.82, 0.37)"
$ws.Range("G3:H3").ClearContents()

# Row 4 (ss3)
$ws.Range("C4").Value2 = "-0.53 (-1.06, 0)"
$ws.Range("E4").Value2 = "-0.12 (-0.65, 0.43)"
$ws.Range("G4").Value2 = "-0.32 (-0.87, 0.26)"
$ws.Range("H4").Value2 = "0.26"

# Row 5 (ss4)
$ws.Range("E5").Value2 = "-0.24 (-0.79, 0.32)"
$ws.Range("G5:H5").ClearContents()

# Row 6 (ss5)
$ws.Range("E6").Value2 = "-0.84 (-1.57, -0.09)"
$ws.Range("G6").Value2 = "-0.93 (-1.69, -0.11)"
$ws.Range("H6").Value2 = "0.02"

# Rows 7-8 now hold aic/bic (shifted up); clear stray "final" column values
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# Row 9 used to be pr_chisq; remove it entirely so the sheet ends at row 8
$ws.Range("A9:H9").ClearContents()
